# "other smx cases added"
# Adds 6 new question-type columns (BX:CC) to the SanityTC sheet header row,
# plus matching blank-but-styled cells in the two data rows, and updates the
# sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SanityTC")
$ws.Activate()

# --- Copy formatting from the last existing column (BW) into the new columns
# so the new cells pick up the same visual style as their neighbours.
$ws.Range("BW1").Copy()
$ws.Range("BX1:CC1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("BW2").Copy()
$ws.Range("BX2:CC2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("BW3").Copy()
$ws.Range("BX3:CC3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- New header values (row 1) ---
$ws.Range("BX1").Value = "horizontalradiobutton"
$ws.Range("BY1").Value = "numericallocations"
$ws.Range("BZ1").Value = "attachments"
$ws.Range("CA1").Value = "ratingradiobutton"
$ws.Range("CB1").Value = "ratingdropdownbutton"
$ws.Range("CC1").Value = "listbox"

# Rows 2 and 3 stay blank in the new columns (same as the rest of the row).

# --- Update the active selection / view to match the edited area ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 62
$win.ScrollRow = 1
$ws.Range("BZ7").Select()
